# Regenerate the K (strikeouts) column values in column G (rows 2-67)
# replacing the old "Strike#" derived values with the recomputed s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sVals = @{
    2  = 1;  3  = 0;  4  = 0;  5  = 2;  6  = 1;  7  = 0;  8  = 1;  9  = 1;  10 = 0;
    11 = 3;  12 = 1;  13 = 1;  14 = 1;  15 = 0;  16 = 1;  17 = 1;  18 = 1;  19 = 0;  20 = 1;
    21 = 2;  22 = 0;  23 = 1;  24 = 0;  25 = 0;  26 = 0;  27 = 2;  28 = 0;  29 = 0;  30 = 1;
    31 = 2;  32 = 0;  33 = 0;  34 = 0;  35 = 1;  36 = 0;  37 = 0;  38 = 0;  39 = 1;  40 = 1;
    41 = 3;  42 = 2;  43 = 0;  44 = 0;  45 = 2;  46 = 0;  47 = 1;  48 = 1;  49 = 2;  50 = 1;
    51 = 1;  52 = 0;  53 = 1;  54 = 1;  55 = 0;  56 = 1;  57 = 0;  58 = 0;  59 = 0;  60 = 0;
    61 = 0;  62 = 1;  63 = 1;  64 = 0;  65 = 2;  66 = 0;  67 = 1
}

foreach ($row in $sVals.Keys) {
    $ws.Cells.Item($row, 7).Value = $sVals[$row]
}
